$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 575.1429000000001
$ws.Range("I2").Value = 482.33334
$ws.Range("K2").Value = 482.33334
$ws.Range("M2").Value = -369.33334
$ws.Range("H11").Value = 1021.2
$ws.Range("I11").Value = 1021.2
$ws.Range("K11").Value = 1021.2
$ws.Range("M11").Value = -881.2
$ws.Range("H18").Value = 1744.8148
$ws.Range("I18").Value = 487
$ws.Range("K18").Value = 487
$ws.Range("M18").Value = -203
$ws.Range("H19").Value = 1568.1052
$ws.Range("I19").Value = 800.7778
$ws.Range("K19").Value = 800.7778
$ws.Range("M19").Value = -625.7778
$ws.Range("H21").Value = 3225
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 3225
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H40").Value = 21433604
$ws.Range("I40").Value = 6581.857
$ws.Range("J40").Value = 42860628
$ws.Range("K40").Value = 6581.857
$ws.Range("L40").Value = 42860628
$ws.Range("M40").Value = -6406.857
$ws.Range("N40").Value = -42860978
$ws.Range("H53").Value = 1213.4
$ws.Range("I53").Value = 1129.909
$ws.Range("K53").Value = 1129.909
$ws.Range("M53").Value = -492.9090000000001
$ws.Range("H62").Value = 9812466
$ws.Range("J62").Value = 12497
$ws.Range("L62").Value = 12497
$ws.Range("N62").Value = -13745
$ws.Range("H65").Value = 9812466
$ws.Range("J65").Value = 12497
$ws.Range("L65").Value = 62485
$ws.Range("N65").Value = -68725
$ws.Range("H70").Value = 6022.486
$ws.Range("I70").Value = 1371.45
$ws.Range("J70").Value = 12223.866
$ws.Range("K70").Value = 4114.35
$ws.Range("L70").Value = 36671.598
$ws.Range("M70").Value = -3844.35
$ws.Range("N70").Value = -37211.598
$ws.Range("H73").Value = 6022.486
$ws.Range("I73").Value = 1371.45
$ws.Range("J73").Value = 12223.866
$ws.Range("K73").Value = 4114.35
$ws.Range("L73").Value = 36671.598
$ws.Range("M73").Value = -3178.35
$ws.Range("N73").Value = -38543.598
$ws.Range("H74").Value = 12026.429
$ws.Range("I74").Value = 12026.429
$ws.Range("K74").Value = 12026.429
$ws.Range("M74").Value = -11090.429
$ws.Range("H76").Value = 2997
$ws.Range("J76").Value = 2997
$ws.Range("L76").Value = 2997
$ws.Range("N76").Value = -3627
$ws.Range("H77").Value = 12026.429
$ws.Range("I77").Value = 12026.429
$ws.Range("K77").Value = 60132.145
$ws.Range("M77").Value = -55452.145
$ws.Range("H79").Value = 2997
$ws.Range("J79").Value = 2997
$ws.Range("L79").Value = 2997
$ws.Range("N79").Value = -5181
$ws.Range("H80").Value = 615.6667
$ws.Range("I80").Value = 615.6667
$ws.Range("K80").Value = 1847.0001
$ws.Range("M80").Value = -849.0001
$ws.Range("H83").Value = 615.6667
$ws.Range("I83").Value = 615.6667
$ws.Range("K83").Value = 5541.0003
$ws.Range("M83").Value = -549.0002999999997
$ws.Range("H86").Value = 4265.125
$ws.Range("I86").Value = 3974.8
$ws.Range("J86").Value = 4749
$ws.Range("K86").Value = 3974.8
$ws.Range("L86").Value = 4749
$ws.Range("M86").Value = -2851.8
$ws.Range("N86").Value = -6995
$ws.Range("H89").Value = 4265.125
$ws.Range("I89").Value = 3974.8
$ws.Range("J89").Value = 4749
$ws.Range("K89").Value = 19874
$ws.Range("L89").Value = 23745
$ws.Range("M89").Value = -14258
$ws.Range("N89").Value = -34977
$ws.Range("H96").Value = 974.9
$ws.Range("I96").Value = 721.2857
$ws.Range("J96").Value = 1566.6666
$ws.Range("K96").Value = 2163.8571
$ws.Range("L96").Value = 4699.9998
$ws.Range("M96").Value = -790.8571000000002
$ws.Range("N96").Value = -7445.9998
$ws.Range("H97").Value = 333339170
$ws.Range("J97").Value = 333339170
$ws.Range("L97").Value = 1000017510
$ws.Range("N97").Value = -1000018502
$ws.Range("H106").Value = 3370.625
$ws.Range("I106").Value = 3030.7856
$ws.Range("K106").Value = 3030.7856
$ws.Range("M106").Value = -2399.7856
$ws.Range("H113").Value = 2987.2964
$ws.Range("I113").Value = 3450.158
$ws.Range("K113").Value = 3450.158
$ws.Range("M113").Value = -196.1579999999999
$ws.Range("H115").Value = 2851738.5
$ws.Range("I115").Value = 3369745.5
$ws.Range("K115").Value = 10109236.5
$ws.Range("M115").Value = -10107669.5
$ws.Range("H121").Value = 4049.5
$ws.Range("J121").Value = 4049.5
$ws.Range("L121").Value = 12148.5
$ws.Range("N121").Value = -15642.5
$ws.Range("H132").Value = 3048.72
$ws.Range("I132").Value = 2605.5908
$ws.Range("K132").Value = 7816.7724
$ws.Range("M132").Value = -5286.7724
$ws.Range("H135").Value = 1554.7059
$ws.Range("I135").Value = 1113.75
$ws.Range("J135").Value = 2613
$ws.Range("K135").Value = 10023.75
$ws.Range("L135").Value = 23517
$ws.Range("M135").Value = -7488.75
$ws.Range("N135").Value = -28587
$ws.Range("H137").Value = 2299286
$ws.Range("I137").Value = 3262338.8
$ws.Range("K137").Value = 9787016.399999999
$ws.Range("M137").Value = -9784466.399999999
$ws.Range("H138").Value = 3557.0952
$ws.Range("I138").Value = 1857.4286
$ws.Range("J138").Value = 4406.9287
$ws.Range("K138").Value = 5572.2858
$ws.Range("L138").Value = 13220.7861
$ws.Range("M138").Value = -432.2857999999997
$ws.Range("N138").Value = -23500.7861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1978
$ws.Range("I2").Value = 1951.5385
$ws.Range("K2").Value = 1951.5385
$ws.Range("M2").Value = -1838.5385
$ws.Range("H21").Value = 1122.6
$ws.Range("I21").Value = 900
$ws.Range("J21").Value = 1271
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1271
$ws.Range("M21").Value = -526
$ws.Range("N21").Value = -2019
$ws.Range("H32").Value = 1557.1915
$ws.Range("I32").Value = 1579.8556
$ws.Range("K32").Value = 1579.8556
$ws.Range("M32").Value = -1292.8556
$ws.Range("H36").Value = 3549.8
$ws.Range("J36").Value = 1916.3334
$ws.Range("L36").Value = 1916.3334
$ws.Range("N36").Value = -2608.3334
$ws.Range("H45").Value = 1971.125
$ws.Range("I45").Value = 1916.6154
$ws.Range("J45").Value = 2207.3333
$ws.Range("K45").Value = 1916.6154
$ws.Range("L45").Value = 2207.3333
$ws.Range("M45").Value = -1539.6154
$ws.Range("N45").Value = -2961.3333
$ws.Range("H61").Value = 5794.727
$ws.Range("I61").Value = 2515.3333
$ws.Range("K61").Value = 2515.3333
$ws.Range("M61").Value = -2303.3333
$ws.Range("H74").Value = 125545.375
$ws.Range("I74").Value = 153386.69
$ws.Range("K74").Value = 153386.69
$ws.Range("M74").Value = -152512.69
$ws.Range("H76").Value = 33808.668
$ws.Range("I76").Value = 31639
$ws.Range("K76").Value = 31639
$ws.Range("M76").Value = -31301
$ws.Range("H77").Value = 125545.375
$ws.Range("I77").Value = 153386.69
$ws.Range("K77").Value = 766933.45
$ws.Range("M77").Value = -762565.45
$ws.Range("H79").Value = 33808.668
$ws.Range("I79").Value = 31639
$ws.Range("K79").Value = 31639
$ws.Range("M79").Value = -30469
$ws.Range("H102").Value = 1675.8823
$ws.Range("I102").Value = 1542.1538
$ws.Range("J102").Value = 2110.5
$ws.Range("K102").Value = 1542.1538
$ws.Range("L102").Value = 2110.5
$ws.Range("M102").Value = 79.84619999999995
$ws.Range("N102").Value = -5354.5
$ws.Range("H116").Value = 1978
$ws.Range("I116").Value = 1951.5385
$ws.Range("K116").Value = 1951.5385
$ws.Range("M116").Value = 342.4614999999999
$ws.Range("H117").Value = 66999.89999999999
$ws.Range("J117").Value = 66999.89999999999
$ws.Range("L117").Value = 66999.89999999999
$ws.Range("N117").Value = -76177.89999999999
$ws.Range("H122").Value = 2196.6924
$ws.Range("I122").Value = 1990.9736
$ws.Range("K122").Value = 5972.9208
$ws.Range("M122").Value = -3522.9208
$ws.Range("H136").Value = 5794.727
$ws.Range("I136").Value = 2515.3333
$ws.Range("K136").Value = 7545.999899999999
$ws.Range("M136").Value = -4995.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1978
$ws.Range("I3").Value = 1951.5385
$ws.Range("K3").Value = 1951.5385
$ws.Range("M3").Value = -1837.5385
$ws.Range("H20").Value = 3164.111
$ws.Range("I20").Value = 3410.6667
$ws.Range("J20").Value = 1931.3334
$ws.Range("K20").Value = 3410.6667
$ws.Range("L20").Value = 1931.3334
$ws.Range("M20").Value = -3163.6667
$ws.Range("N20").Value = -2425.3334
$ws.Range("H94").Value = 798.14703
$ws.Range("I94").Value = 668.5
$ws.Range("J94").Value = 1403.1666
$ws.Range("K94").Value = 668.5
$ws.Range("L94").Value = 1403.1666
$ws.Range("M94").Value = -217.5
$ws.Range("N94").Value = -2305.1666
$ws.Range("H99").Value = 5272.353
$ws.Range("I99").Value = 6175
$ws.Range("J99").Value = 4994.615
$ws.Range("K99").Value = 6175
$ws.Range("L99").Value = 4994.615
$ws.Range("M99").Value = -4677
$ws.Range("N99").Value = -7990.615
$ws.Range("H105").Value = 2062.2104
$ws.Range("I105").Value = 2134.8333
$ws.Range("J105").Value = 1937.7142
$ws.Range("K105").Value = 2134.8333
$ws.Range("L105").Value = 1937.7142
$ws.Range("M105").Value = -387.8332999999998
$ws.Range("N105").Value = -5431.7142
$ws.Range("H107").Value = 1438.4
$ws.Range("I107").Value = 799
$ws.Range("J107").Value = 3996
$ws.Range("K107").Value = 799
$ws.Range("L107").Value = 3996
$ws.Range("M107").Value = 1121
$ws.Range("N107").Value = -7836
$ws.Range("H134").Value = 7857.625
$ws.Range("I134").Value = 7857.625
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 23572.875
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -21037.875
$ws.Range("N134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 730.4545000000001
$ws.Range("I7").Value = 566.2381
$ws.Range("J7").Value = 1017.8333
$ws.Range("K7").Value = 566.2381
$ws.Range("L7").Value = 1017.8333
$ws.Range("M7").Value = -453.2381
$ws.Range("N7").Value = -1243.8333
$ws.Range("H16").Value = 1350.56
$ws.Range("I16").Value = 776.6842
$ws.Range("J16").Value = 3167.8333
$ws.Range("K16").Value = 776.6842
$ws.Range("L16").Value = 3167.8333
$ws.Range("M16").Value = -489.6842
$ws.Range("N16").Value = -3741.8333
$ws.Range("H22").Value = 320.3611
$ws.Range("I22").Value = 321.90625
$ws.Range("J22").Value = 308
$ws.Range("K22").Value = 321.90625
$ws.Range("L22").Value = 308
$ws.Range("M22").Value = 28.09375
$ws.Range("N22").Value = -1008
$ws.Range("H31").Value = 272043.84
$ws.Range("I31").Value = 668667.75
$ws.Range("J31").Value = 1618.4546
$ws.Range("K31").Value = 668667.75
$ws.Range("L31").Value = 1618.4546
$ws.Range("M31").Value = -668372.75
$ws.Range("N31").Value = -2208.4546
$ws.Range("H34").Value = 272043.84
$ws.Range("I34").Value = 668667.75
$ws.Range("J34").Value = 1618.4546
$ws.Range("K34").Value = 668667.75
$ws.Range("L34").Value = 1618.4546
$ws.Range("M34").Value = -668465.75
$ws.Range("N34").Value = -2022.4546
$ws.Range("H36").Value = 10000
$ws.Range("J36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("N36").Value = -10776
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10320
$ws.Range("H58").Value = 1927.3334
$ws.Range("I58").Value = 1893.4546
$ws.Range("K58").Value = 1893.4546
$ws.Range("M58").Value = -1690.4546
$ws.Range("H113").Value = 1350.56
$ws.Range("I113").Value = 776.6842
$ws.Range("J113").Value = 3167.8333
$ws.Range("K113").Value = 776.6842
$ws.Range("L113").Value = 3167.8333
$ws.Range("M113").Value = 1393.3158
$ws.Range("N113").Value = -7507.8333
$ws.Range("H132").Value = 6370.757
$ws.Range("I132").Value = 6335.0625
$ws.Range("K132").Value = 19005.1875
$ws.Range("M132").Value = -16475.1875
$ws.Range("H134").Value = 2660.2
$ws.Range("I134").Value = 2409.1428
$ws.Range("K134").Value = 7227.428400000001
$ws.Range("M134").Value = -4692.428400000001
$ws.Range("H136").Value = 1927.3334
$ws.Range("I136").Value = 1893.4546
$ws.Range("K136").Value = 5680.3638
$ws.Range("M136").Value = -3130.3638

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 597.46155
$ws.Range("I34").Value = 110.4
$ws.Range("J34").Value = 2221
$ws.Range("K34").Value = 331.2
$ws.Range("L34").Value = 6663
$ws.Range("M34").Value = -247.2
$ws.Range("N34").Value = -6831
$ws.Range("H39").Value = 1000.38464
$ws.Range("I39").Value = 834.5833
$ws.Range("K39").Value = 2503.7499
$ws.Range("M39").Value = -2209.7499
$ws.Range("H55").Value = 8996.429
$ws.Range("J55").Value = 8996.429
$ws.Range("L55").Value = 26989.287
$ws.Range("N55").Value = -27343.287
$ws.Range("H68").Value = 2779514.5
$ws.Range("I68").Value = 20834134
$ws.Range("J68").Value = 1880.6154
$ws.Range("K68").Value = 62502402
$ws.Range("L68").Value = 5641.8462
$ws.Range("M68").Value = -62501591
$ws.Range("N68").Value = -7263.8462
$ws.Range("H71").Value = 2779514.5
$ws.Range("I71").Value = 20834134
$ws.Range("J71").Value = 1880.6154
$ws.Range("K71").Value = 187507206
$ws.Range("L71").Value = 16925.5386
$ws.Range("M71").Value = -187503150
$ws.Range("N71").Value = -25037.5386
$ws.Range("H97").Value = 349.75
$ws.Range("J97").Value = 99.333336
$ws.Range("L97").Value = 298.000008
$ws.Range("N97").Value = -1290.000008
$ws.Range("H107").Value = 1003.4167
$ws.Range("I107").Value = 376
$ws.Range("J107").Value = 1317.125
$ws.Range("K107").Value = 1128
$ws.Range("L107").Value = 3951.375
$ws.Range("M107").Value = 792
$ws.Range("N107").Value = -7791.375
$ws.Range("H109").Value = 201246.44
$ws.Range("I109").Value = 201246.44
$ws.Range("K109").Value = 603739.3200000001
$ws.Range("M109").Value = -602699.3200000001
$ws.Range("H114").Value = 2180.2222
$ws.Range("I114").Value = 344.6
$ws.Range("J114").Value = 4474.75
$ws.Range("K114").Value = 1033.8
$ws.Range("L114").Value = 13424.25
$ws.Range("M114").Value = 2220.2
$ws.Range("N114").Value = -19932.25
$ws.Range("H117").Value = 257.75
$ws.Range("I117").Value = 257.75
$ws.Range("K117").Value = 773.25
$ws.Range("M117").Value = 2668.75
$ws.Range("H131").Value = 5557176
$ws.Range("I131").Value = 50001148
$ws.Range("J131").Value = 1679.3125
$ws.Range("K131").Value = 150003444
$ws.Range("L131").Value = 5037.9375
$ws.Range("M131").Value = -149998404
$ws.Range("N131").Value = -15117.9375
$ws.Range("H132").Value = 9882.214
$ws.Range("I132").Value = 11004.25
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 99038.25
$ws.Range("L132").Value = 28350
$ws.Range("M132").Value = -96508.25
$ws.Range("N132").Value = -33410

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H70").Value = 10407.895
$ws.Range("I70").Value = 8337.700000000001
$ws.Range("J70").Value = 12708.111
$ws.Range("K70").Value = 8337.700000000001
$ws.Range("L70").Value = 12708.111
$ws.Range("M70").Value = -8067.700000000001
$ws.Range("N70").Value = -13248.111
$ws.Range("H73").Value = 10407.895
$ws.Range("I73").Value = 8337.700000000001
$ws.Range("J73").Value = 12708.111
$ws.Range("K73").Value = 8337.700000000001
$ws.Range("L73").Value = 12708.111
$ws.Range("M73").Value = -7401.700000000001
$ws.Range("N73").Value = -14580.111
$ws.Range("H97").Value = 2478.5
$ws.Range("I97").Value = 2798.6667
$ws.Range("J97").Value = 1998.25
$ws.Range("K97").Value = 2798.6667
$ws.Range("L97").Value = 1998.25
$ws.Range("M97").Value = -2302.6667
$ws.Range("N97").Value = -2990.25
$ws.Range("H102").Value = 3489.4285
$ws.Range("I102").Value = 2708.4
$ws.Range("J102").Value = 9998
$ws.Range("K102").Value = 2708.4
$ws.Range("L102").Value = 9998
$ws.Range("M102").Value = -1086.4
$ws.Range("N102").Value = -13242
$ws.Range("H113").Value = 9582.833000000001
$ws.Range("I113").Value = 14332.333
$ws.Range("J113").Value = 4833.3335
$ws.Range("K113").Value = 14332.333
$ws.Range("L113").Value = 4833.3335
$ws.Range("M113").Value = -12162.333
$ws.Range("N113").Value = -9173.333500000001
$ws.Range("H122").Value = 2801.9375
$ws.Range("I122").Value = 2478.1667
$ws.Range("K122").Value = 7434.500100000001
$ws.Range("M122").Value = -4984.500100000001
$ws.Range("H123").Value = 39460.8
$ws.Range("J123").Value = 39460.8
$ws.Range("L123").Value = 39460.8
$ws.Range("N123").Value = -44360.8
$ws.Range("H132").Value = 44055.25
$ws.Range("I132").Value = 54659.58
$ws.Range("J132").Value = 3758.8
$ws.Range("K132").Value = 163978.74
$ws.Range("L132").Value = 11276.4
$ws.Range("M132").Value = -161448.74
$ws.Range("N132").Value = -16336.4
$ws.Range("H133").Value = 116236.25
$ws.Range("J133").Value = 116236.25
$ws.Range("L133").Value = 116236.25
$ws.Range("N133").Value = -126356.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5365.095
$ws.Range("I7").Value = 5279.125
$ws.Range("J7").Value = 5640.2
$ws.Range("K7").Value = 5279.125
$ws.Range("L7").Value = 5640.2
$ws.Range("M7").Value = -5167.125
$ws.Range("N7").Value = -5864.2
$ws.Range("H22").Value = 697.2222
$ws.Range("I22").Value = 681.1667
$ws.Range("J22").Value = 729.3333
$ws.Range("K22").Value = 681.1667
$ws.Range("L22").Value = 729.3333
$ws.Range("M22").Value = -386.1667
$ws.Range("N22").Value = -1319.3333
$ws.Range("H27").Value = 697.2222
$ws.Range("I27").Value = 681.1667
$ws.Range("J27").Value = 729.3333
$ws.Range("K27").Value = 681.1667
$ws.Range("L27").Value = 729.3333
$ws.Range("M27").Value = -574.1667
$ws.Range("N27").Value = -943.3333
$ws.Range("H46").Value = 2875.2
$ws.Range("J46").Value = 3609.889
$ws.Range("L46").Value = 3609.889
$ws.Range("N46").Value = -3985.889
$ws.Range("H61").Value = 1969.625
$ws.Range("I61").Value = 1969.625
$ws.Range("K61").Value = 1969.625
$ws.Range("M61").Value = -1767.625
$ws.Range("H82").Value = 1829.3077
$ws.Range("I82").Value = 1751.2106
$ws.Range("K82").Value = 1751.2106
$ws.Range("M82").Value = -1390.2106
$ws.Range("H85").Value = 1829.3077
$ws.Range("I85").Value = 1751.2106
$ws.Range("K85").Value = 1751.2106
$ws.Range("M85").Value = -503.2106000000001
$ws.Range("H100").Value = 26666.666
$ws.Range("I100").Value = 30000
$ws.Range("K100").Value = 30000
$ws.Range("M100").Value = -29459
$ws.Range("H113").Value = 1969.625
$ws.Range("I113").Value = 1969.625
$ws.Range("K113").Value = 1969.625
$ws.Range("M113").Value = 200.375
$ws.Range("H122").Value = 4178.727
$ws.Range("I122").Value = 4331.778
$ws.Range("K122").Value = 12995.334
$ws.Range("M122").Value = -10545.334
$ws.Range("H126").Value = 5365.095
$ws.Range("I126").Value = 5279.125
$ws.Range("J126").Value = 5640.2
$ws.Range("K126").Value = 15837.375
$ws.Range("L126").Value = 16920.6
$ws.Range("M126").Value = -13367.375
$ws.Range("N126").Value = -21860.6
$ws.Range("H132").Value = 4160.25
$ws.Range("I132").Value = 4123.5
$ws.Range("K132").Value = 12370.5
$ws.Range("M132").Value = -9840.5
$ws.Range("H136").Value = 3014.5806
$ws.Range("I136").Value = 3014.5217
$ws.Range("K136").Value = 9043.5651
$ws.Range("M136").Value = -6493.5651

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 5500
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 5500
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 5500
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -6446
$ws.Range("H57").Value = 49999
$ws.Range("I57").Value = 49999
$ws.Range("K57").Value = 49999
$ws.Range("M57").Value = -49245
$ws.Range("H62").Value = 9886.223
$ws.Range("I62").Value = 7499
$ws.Range("J62").Value = 10568.286
$ws.Range("K62").Value = 7499
$ws.Range("L62").Value = 10568.286
$ws.Range("M62").Value = -6875
$ws.Range("N62").Value = -11816.286
$ws.Range("H63").Value = 2369.5
$ws.Range("J63").Value = 2369.5
$ws.Range("L63").Value = 2369.5
$ws.Range("N63").Value = -3617.5
$ws.Range("H65").Value = 9886.223
$ws.Range("I65").Value = 7499
$ws.Range("J65").Value = 10568.286
$ws.Range("K65").Value = 37495
$ws.Range("L65").Value = 52841.43
$ws.Range("M65").Value = -34375
$ws.Range("N65").Value = -59081.43
$ws.Range("H66").Value = 2369.5
$ws.Range("J66").Value = 2369.5
$ws.Range("L66").Value = 7108.5
$ws.Range("N66").Value = -13348.5
$ws.Range("H122").Value = 2073.2856
$ws.Range("I122").Value = 2073.2856
$ws.Range("K122").Value = 6219.8568
$ws.Range("M122").Value = -3769.8568
$ws.Range("H126").Value = 3178.5
$ws.Range("I126").Value = 2989.111
$ws.Range("J126").Value = 3519.4
$ws.Range("K126").Value = 8967.332999999999
$ws.Range("L126").Value = 10558.2
$ws.Range("M126").Value = -6497.332999999999
$ws.Range("N126").Value = -15498.2
$ws.Range("H132").Value = 2179.9207
$ws.Range("I132").Value = 1827.7838
$ws.Range("J132").Value = 2681.0386
$ws.Range("K132").Value = 5483.3514
$ws.Range("L132").Value = 8043.1158
$ws.Range("M132").Value = -2953.3514
$ws.Range("N132").Value = -13103.1158
$ws.Range("H136").Value = 419506.84
$ws.Range("I136").Value = 457439.3
$ws.Range("K136").Value = 1372317.9
$ws.Range("M136").Value = -1369767.9
